$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the run "The responsible grammar rules are:" into two runs
#           "The responsible grammar r" + "ules are:" (no text change, no
#           bookmark left behind).
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("The responsible grammar rules are:")
$splitPos = $idx + ("The responsible grammar r").Length
$rng = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplit1", $rng) | Out-Null
$d.Bookmarks("TempSplit1").Delete()

# ---------------------------------------------------------------------------
# Change 2: fix "1 PP Prep PP" -> "1 PP Prep NP", splitting it into
#           "1 PP Prep N" + "P" with a "_GoBack" bookmark between them.
#           A temporary bookmark is placed just before the run (before its
#           leading <w:br/>) so the text edit cannot coalesce this run back
#           into the preceding " PP" run.
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$idx = $text.IndexOf("1 PP Prep PP")

$blockPos = $idx - 1
$blockRng = $d.Range($blockPos, $blockPos)
$d.Bookmarks.Add("TempBlock2", $blockRng) | Out-Null

$fixRng = $d.Range($idx + 10, $idx + 11)
$fixRng.Text = "N"

$d.Bookmarks("TempBlock2").Delete()

$splitPos2 = $idx + 11
$rng2 = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $rng2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: merge the " PP" run with the " ' still has weight of 1, but the
#           'NP " run that follows it into a single run. We force this by
#           performing a genuine (non no-op) text edit on the second run --
#           first to a placeholder, then back to the original text -- which
#           triggers the engine's run-coalescing with the preceding,
#           identically-formatted run.
# ---------------------------------------------------------------------------
$target3 = " ' still has weight of 1, but the 'NP "
$text = $d.Content.Text
$idx3 = $text.IndexOf($target3)
$fixRng3 = $d.Range($idx3, $idx3 + $target3.Length)
$fixRng3.Text = "TempPlaceholder3"

$text = $d.Content.Text
$idx3b = $text.IndexOf("TempPlaceholder3")
$fixRng3b = $d.Range($idx3b, $idx3b + "TempPlaceholder3".Length)
$fixRng3b.Text = $target3

# ---------------------------------------------------------------------------
# Change 4: remove the old "_GoBack" bookmark that used to sit at the end of
#           the last paragraph (it moved to the middle of "1 PP Prep NP").
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()
$rngEnd = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $rngEnd) | Out-Null

Write-Host "done"
